$wb = $excel.ActiveWorkbook

# Rename the simulation sheets to reflect the wd= parameter values used
# in this run (011220 results).
$wb.Worksheets.Item(1).Name = "wd=0"
$wb.Worksheets.Item(2).Name = "wd=0.000001"
$wb.Worksheets.Item(3).Name = "wd=0.00001"
$wb.Worksheets.Item(4).Name = "wd=0.0001"
$wb.Worksheets.Item(5).Name = "wd=0.001"

# Make the first sheet ("wd=0") the active/selected tab instead of the
# second one.
$wb.Worksheets.Item(1).Select()
